# Results from December_05,_2020--22:02:44 run
# Applies the daily COVID disparities data refresh (2020-12-04 -> 2020-12-05)
# to the rows whose upstream source was re-scraped in this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric cell updates (Total Cases/Deaths, counts, pct's, dates as serials) ---
$ws.Range("B2").Value = 44170
$ws.Range("C2").Value = 779975
$ws.Range("D2").Value = 13179
$ws.Range("E2").Value = 72772
$ws.Range("F2").Value = 2687
$ws.Range("G2").Value = 9.33
$ws.Range("H2").Value = 20.39

$ws.Range("C4").Value = 175793
$ws.Range("E4").Value = 6253
$ws.Range("F4").Value = 89
$ws.Range("G4").Value = 5.84
$ws.Range("H4").Value = 3.1
$ws.Range("K4").Value = 106997
$ws.Range("L4").Value = 2875

$ws.Range("B5").Value = 44170
$ws.Range("C5").Value = 388552
$ws.Range("D5").Value = 5516
$ws.Range("E5").Value = 69499
$ws.Range("F5").Value = 1476
$ws.Range("G5").Value = 21.72
$ws.Range("H5").Value = 28.03
$ws.Range("K5").Value = 320018
$ws.Range("L5").Value = 5266

$ws.Range("B12").Value = 44170
$ws.Range("C12").Value = 397522
$ws.Range("D12").Value = 4905
$ws.Range("E12").Value = 56328
$ws.Range("F12").Value = 1017
$ws.Range("G12").Value = 14.17
$ws.Range("H12").Value = 20.73

$ws.Range("B16").Value = 44170
$ws.Range("C16").Value = 169382
$ws.Range("D16").Value = 2620
$ws.Range("E16").Value = 29954
$ws.Range("F16").Value = 428
$ws.Range("G16").Value = 20.03
$ws.Range("H16").Value = 17.4
$ws.Range("K16").Value = 149579
$ws.Range("L16").Value = 2460

$ws.Range("B19").Value = 44170
$ws.Range("C19").Value = 394976
$ws.Range("D19").Value = 9797
$ws.Range("E19").Value = 45609
$ws.Range("F19").Value = 2780
$ws.Range("G19").Value = 11.55
$ws.Range("H19").Value = 28.38

# Row 21 (Oklahoma): B21 stays text "2020-12-05" (see below); rest are numbers.
$ws.Range("C21").Value = 213245
$ws.Range("D21").Value = 1874
$ws.Range("E21").Value = 12773.3755
$ws.Range("F21").Value = 113.0022
$ws.Range("G21").Value = 5.99
$ws.Range("H21").Value = 6.03
$ws.Range("K21").Value = 166331.1
$ws.Range("L21").Value = 1716.0218

$ws.Range("B24").Value = 44170
$ws.Range("C24").Value = 27861
$ws.Range("D24").Value = 141
$ws.Range("E24").Value = 943
$ws.Range("G24").Value = 3.44
$ws.Range("H24").Value = 4.26
$ws.Range("K24").Value = 27429
$ws.Range("L24").Value = 141

$ws.Range("C27").Value = 243384
$ws.Range("D27").Value = 2666
$ws.Range("E27").Value = 8162
$ws.Range("G27").Value = 3.35

# Row 29 (Utah): C29/D29/E29 stay text (see below); B29/H29 are numbers.
$ws.Range("B29").Value = 44170
$ws.Range("H29").Value = 1.17

$ws.Range("B35").Value = 44170
$ws.Range("C35").Value = 138568
$ws.Range("D35").Value = 1194

$ws.Range("B36").Value = 44169
$ws.Range("C36").Value = 1311625
$ws.Range("D36").Value = 19791
$ws.Range("E36").Value = 38841
$ws.Range("F36").Value = 1418
$ws.Range("H36").Value = 7.23
$ws.Range("K36").Value = 947773
$ws.Range("L36").Value = 19602

$ws.Range("B39").Value = 44170
$ws.Range("C39").Value = 106856
$ws.Range("D39").Value = 1738
$ws.Range("E39").Value = 1312

$ws.Range("B42").Value = 44170
$ws.Range("C42").Value = 257347
$ws.Range("D42").Value = 3357
$ws.Range("E42").Value = 6202
$ws.Range("G42").Value = 3.54
$ws.Range("H42").Value = 4.75
$ws.Range("K42").Value = 175279
$ws.Range("L42").Value = 3157

$ws.Range("B49").Value = 44170
$ws.Range("C49").Value = 375019
$ws.Range("D49").Value = 5910
$ws.Range("E49").Value = 25137
$ws.Range("F49").Value = 569
$ws.Range("G49").Value = 6.7
$ws.Range("H49").Value = 9.630000000000001

# --- Text cells that must remain stored as shared-string text (not auto-coerced
# to a date serial / number by Excel's smart-entry). Force text via NumberFormat
# "@" before assigning, then drop back to the "Normal" style so no visible
# number-format change is left behind (matches the unstyled source cells). ---

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "2020-12-05"
$ws.Range("B21").Style = "Normal"

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "212844"
$ws.Range("C29").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "939"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "3078"
$ws.Range("E29").Style = "Normal"

# --- Status-code message for Florida (O40): the pipeline re-ran and the
# temp pdf filename embedded in the tabula error message changed. ---
$ws.Range("O40").Value = "An error occurred. ... CalledProcessError(1, ['java', '-Djava.awt.headless=true', '-Dfile.encoding=UTF8', '-jar', '/Users/poisson/Documents/GitHub/COVID19_tracker_data_extraction/covid19_data_test_003/lib/python3.7/site-packages/tabula/tabula-1.0.3-jar-with-dependencies.jar', '--pages', '3', '--area', '626.8578491210938,80.14600372314453,961.4368286132812,941.5399780273438', '--stream', '/var/folders/cn/4q_jgh710l170_p6btg_ym880000gn/T/a79b4226-8bf5-4467-b1e6-407a997a8bdd.pdf'])"
